$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target table (header + 18 data rows) as it should look after the edit.
$data = @(
    @("Oyuncu Adı", "Pozisyon", "Takım"),
    @("De'Aaron Fox", "PG", "Sacramento Kings"),
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Gradey Dick", "SG,SF", "Toronto Raptors"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Bennedict Mathurin", "SG,SF", "Indiana Pacers"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("Luka Doncic", "PG,SG", "Dallas Mavericks"),
    @("Amen Thompson", "SG,SF", "Houston Rockets"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Caris LeVert", "SG,SF", "Cleveland Cavaliers"),
    @("Nick Richards", "C", "Charlotte Hornets"),
    @("Scottie Barnes", "SG,SF,PF", "Toronto Raptors"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets"),
    @("Ja Morant", "PG", "Memphis Grizzlies")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
}
